# "Enhanced project and Added Test Case Verify Menu Items"
#
# The "INVALID CREDENTIALS" test block used vj2805 / vishalj28051 as its
# first (now obsolete) pair of sample credentials in row 4. Remove that
# whole row so every row below it (the INVALID CREDENTIALS header and its
# three username/password rows) shifts up by one, exactly as Excel does
# when you delete an entire row from the worksheet.

$wb = $excel.ActiveWorkbook
try {
    $ws = $wb.Worksheets.Item("Credentials")
} catch {
    $ws = $wb.ActiveSheet
}

# Delete the entire 4th row (A4:B4 -> "vj2805" / "vishalj28051"),
# shifting everything below it up by one row.
$ws.Rows.Item(4).EntireRow.Delete()

# Excel leaves a trace of row-level formatting on the very last row of the
# sheet after a full-row delete; reproduce that so the sheet's used range
# stretches down to the bottom of the worksheet.
$ws.Rows.Item($ws.Rows.Count).RowHeight = 12.8

# Reflect the resulting selection: after the shift, the last data row is
# now row 9 and the previously-selected B7 cell content now lives at B9.
$ws.Range("B9").Select()
